$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = 'điều 22. điểm học phần 1. điểm đánh giá thành phần và điểm thi kết thúc học phần được chấm theo thang điểm 10 từ 0 đến 10, làm tròn đến một chữ số thập phân. 2. điểm học phần là tổng số điểm của tất cả các điểm đánh giá thành phần của học phần nhân với trọng số tương ứng. điểm học phần được tính theo thang điểm 10 và làm tròn đến một chữ số thập phân. GV phụ trách học phần nhập điểm vào hệ thống quản lý trực tuyến, hệ thống quy đổi sang điểm chữ và điểm số theo thang điểm 4. cách quy đổi điểm được thực hiện theo bảng dưới đây  điểm số theo thang điểm 10  điểm chữ  điểm số theo thang điểm 4  -----------------------------------------------------------------  9.0 - 10.0  A  4.0   8.0 - 8.9  B+  3.5   7.0 - 7.9  B  3.0   6.5 - 6.9  C+  2.5   5.5 - 6.4  C  2.0   5.0 - 5.4  D+  1.5   4.0 - 4.9  D  1.0   nhỏ hơn 4.0  F  0.0  3. học phần chỉ được tính tích lũy khi đạt từ điểm D trở lên. 4. điểm học phần sẽ được công bố và ghi nhận với điểm số theo thang điểm 10 và điểm chữ. điểm chữ được quy đổi sang điểm số theo thang điểm 4 để tính ĐTBCHK và ĐTBCTL xem điều 24. 5. các điểm học phần khác a điểm M tương đương điểm R của thông tư 08/2021/TT-BGDĐT dùng để xác nhận học phần SV được miễn học do đã tích lũy được bằng các hình thức khác. điểm M không được tính vào ĐTBCHK và ĐTBCTL. số TC của học phần có điểm M được tính vào tổng số TC tích lũy. để nhận điểm M, SV phải làm đơn có ý kiến của trưởng đơn vị đào tạo xem xét theo từng HK, kèm theo giấy tờ hợp lệ. b điểm I chỉ áp dụng cho đánh giá kết thúc học phần. dành cho các trường hợp SV đã dự học, dự các lần kiểm tra giữa HK, đã thực hiện các hoạt động liên quan đến học phần như thí nghiệm, thực hành nhưng vì lý do bất khả kháng như ốm đau, tai nạn,... đã vắng mặt trong buổi thi kết thúc học phần và được GV phụ trách học phần chấp thuận cho bổ sung điểm. điểm I không được tính vào ĐTBCHK ở HK đó. để nhận được điểm I, SV phải làm đơn kèm theo hồ sơ hợp lệ nộp cho GV giảng dạy học phần xem xét và trình trưởng đơn vị đào tạo duyệt. thời hạn bổ sung điểm của học phần do GV quy định nhưng không quá 1 năm kể từ ngày thi lần trước. qua thời hạn trên, nếu SV không hoàn thành thì học phần sẽ nhận điểm F. c điểm W dành cho các học phần mà SV được phép rút theo quy định xem điều 15. điểm W không tính vào ĐTBCHK và ĐTBCTL. '

$ws.Range("A19").Value = $newText

$ws.Rows.Item(23).RowHeight = 188.5
$ws.Rows.Item(24).RowHeight = 87

$ws.Range("A20").Select() | Out-Null
